$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 4 new rows (19-22) that repeat the book list found in rows 15-18,
# continuing the numbering in column A.
$data = @(
    @(18, "The Passionate Programmer", "Chad Fowler", 16),
    @(19, "Software Craftmanship", "Pete McBreen", 26),
    @(20, "The Art of Agile Development", "James Shore", 32),
    @(21, "Continuous Delivery", "Jez Humble", 41)
)

$row = 19
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $row = $row + 1
}
